$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "Samsung Galaxy Tab 10.1" from D1 into C1 (replacing "Canon EOS 5D"),
# then clear out the now-unused D1 cell.
$ws.Range("C1").Value2 = $ws.Range("D1").Value2
$ws.Range("D1").ClearContents()

# Resize column A back to its auto-fit width for the shorter "MacBook"/"iPhone"
# values, and give column C an explicit auto-fit width for its new content
# (column B keeps the sheet's default width, as before).
$ws.Columns("A").ColumnWidth = 8
$ws.Columns("C").ColumnWidth = 21

# Update the active selection to match the post-edit state.
$ws.Range("D5").Select() | Out-Null
